$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 / column B ("B11") held the rule name "R40". Replace it with "1".
$ws.Range("B11").Value = "1"

$wb.Save()
